$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A32").Value = 45946
$ws.Range("B32").Value = "四方坪站"
$ws.Range("C32").Value = 8637.61
$ws.Range("D32").Value = 6843.04
$ws.Range("E32").Value = 2985.93
$ws.Range("F32").Value = 374

$ws.Range("A33").Value = 45946
$ws.Range("B33").Value = "高岭站"
$ws.Range("C33").Value = 3972.73
$ws.Range("D33").Value = 3273.12
$ws.Range("E33").Value = 1042.71
$ws.Range("F33").Value = 146

$ws.Range("A32:A33").NumberFormat = "[$-F800]dddd\,\ mmmm\ dd\,\ yyyy"
$ws.Range("C32:E33").NumberFormat = "0.00_);[Red]\(0.00\)"
$ws.Range("F32:F33").NumberFormat = "0_);[Red]\(0\)"

$ws.Range("H34").Select()
